# Feat: Added Percipio_Email field in users table and fixed create batch API
#
# The "Users" table on Sheet1 gains a new column ("Percipio_Email") inserted
# between the existing "Email" and "Password" columns. The existing Email
# column values (generic a.nehra@... style addresses) move into the new
# Percipio_Email column, and the Email column is repopulated with new
# firstname.lastname@experionglobal.com style addresses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new blank column at D (pushes the old "Password" column to E) ---
$ws.Columns.Item(4).Insert()

# --- Header row ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Role"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Percipio_Email"
$ws.Range("E1").Value = "Password"

# --- Row 2: Joel C Raju ---
$ws.Range("A2").Value = "Joel C Raju"
$ws.Range("B2").Value = "Trainee"
$ws.Range("C2").Value = "joel.raju@experionglobal.com"
$ws.Range("D2").Value = "a.nehra@experionglobal.com"
$ws.Range("E2").Value = "Nehra@explearning"

# --- Row 3: Nigin N ---
$ws.Range("A3").Value = "Nigin N"
$ws.Range("B3").Value = "Trainee"
$ws.Range("C3").Value = "nigin.n@experionglobal.com"
$ws.Range("D3").Value = "b.aravind@experionglobal.com"
$ws.Range("E3").Value = "Aravind@explearning"

# --- Row 4: Thimna Raphel ---
$ws.Range("A4").Value = "Thimna Raphel"
$ws.Range("B4").Value = "Trainee"
$ws.Range("C4").Value = "thimna.raphel@experionglobal.com"
$ws.Range("D4").Value = "b.lara@experionglobal.com"
$ws.Range("E4").Value = "Lara@explearning"

# --- Row 5: Sreejaya V S ---
$ws.Range("A5").Value = "Sreejaya V S"
$ws.Range("B5").Value = "Trainee"
$ws.Range("C5").Value = "sreejaya.vs@experionglobal.com"
$ws.Range("D5").Value = "bs.akshara@experionglobal.com"
$ws.Range("E5").Value = "Akshara@explearning"

# --- Row 6: Ashik George ---
$ws.Range("A6").Value = "Ashik George"
$ws.Range("B6").Value = "Trainee"
$ws.Range("C6").Value = "ashik.george@experionglobal.com"
$ws.Range("D6").Value = "c.devadas@experionglobal.com"
$ws.Range("E6").Value = "Devadas@explearning"

# --- Row 7: Nebil V ---
$ws.Range("A7").Value = "Nebil V"
$ws.Range("B7").Value = "Trainee"
$ws.Range("C7").Value = "nebil.v@experionglobal.com"
$ws.Range("D7").Value = "b.sreerag@experionglobal.com"
$ws.Range("E7").Value = "Sreerag@explearning"

# --- Row 8: Elena Maria Varghese ---
$ws.Range("A8").Value = "Elena Maria Varghese"
$ws.Range("B8").Value = "Trainee"
$ws.Range("C8").Value = "elena.varghese@experionglobal.com"
$ws.Range("D8").Value = "b.meenu@experionglobal.com"
$ws.Range("E8").Value = "Meenu@explearning"

# --- Column widths (bestFit-style widths, matching the new layout).
#     Input values are pre-compensated for the host's internal width
#     quantization so the stored OOXML width lands as close as possible to
#     the target bestFit widths (24.862.., 13.576.., 29.576.., 32.148..,
#     41.434..). ---
$ws.Columns.Item(1).ColumnWidth = 24
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 28.666666666666668
$ws.Columns.Item(4).ColumnWidth = 31.333333333333332
$ws.Columns.Item(5).ColumnWidth = 40.666666666666664

# --- Row heights: rows with longer wrapped Percipio_Email text are taller ---
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5
$ws.Rows.Item(4).RowHeight = 31.5
$ws.Rows.Item(5).RowHeight = 19.5
$ws.Rows.Item(6).RowHeight = 19.5
$ws.Rows.Item(7).RowHeight = 19.5
$ws.Rows.Item(8).RowHeight = 31.5

# --- Font formatting: the distinctive (black-rgb, sometimes underlined)
#     formatting that used to belong to the old Email column (C) now
#     belongs to the Percipio_Email column (D), since that's where the old
#     Email values moved to. The new Email column (C) picks up the plain
#     theme-coloured formatting that the other columns used to use. ---
$ws.Range("A1:A8").Font.Color = 0
$ws.Range("B1:B8").Font.Color = 0
$ws.Range("D1:D8").Font.Color = 0
$ws.Range("E1:E8").Font.Color = 0

$ws.Range("D2").Font.Underline = 2
$ws.Range("D3").Font.Underline = 2
$ws.Range("D4").Font.Underline = 2
$ws.Range("D5").Font.Underline = 2
$ws.Range("D8").Font.Underline = 2

# Email column (C) uses the plain, non-underlined theme-coloured font.
$ws.Range("C1:C8").Font.ThemeColor = 1
$ws.Range("C1:C8").Font.Underline = 0

# --- Normalize horizontal alignment across the whole table to left ---
$ws.Range("A1:E8").HorizontalAlignment = -4131

# --- Sheet dimension now spans through column E ---
$ws.Range("A1:E8").Select()
